# refactor: consistently name Dialogs and put in same folder
#
# The Todo item "change names of dialog-only forms to XDialog" (row 13 on the
# "Active" sheet) is complete. Move it from the "Active" sheet to the
# "Inactive" sheet (as the new first data row), mark its Status as "Done",
# and record the date it was finished ("Done" column) as 8/24/2018.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")

# --- capture the data currently sitting in row 13 of "Active" ---
$taskId       = $active.Range("A13").Text
$taskTitle    = $active.Range("B13").Text
$taskCategory = $active.Range("D13").Text
$taskCreated  = $active.Range("E13").Text

# --- remove that row from "Active"; rows below shift up automatically ---
$active.Rows.Item(13).Delete()

# --- make room for the finished task at the top of "Inactive" ---
$inactive.Rows.Item(2).Insert()

# --- populate the newly inserted row ---
$inactive.Range("A2").Value = $taskId
$inactive.Range("B2").Value = $taskTitle
$inactive.Range("C2").Value = "Done"
$inactive.Range("D2").Value = $taskCategory
$inactive.Range("E2").Value = "'" + $taskCreated
$inactive.Range("F2").Value = "'8/24/2018"
